$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 246-355: rolling weekly price data shift ---
$ws.Range("D246").Value = 44875
$ws.Range("J246").Value = 2400
$ws.Range("K246").Value = 550
$ws.Range("L246").Value = 600
$ws.Range("M246").Value = 575
$ws.Range("P246").Value = 192
$ws.Range("D247").Value = 44875
$ws.Range("K247").Value = 450
$ws.Range("L247").Value = 500
$ws.Range("M247").Value = 475
$ws.Range("P247").Value = 158
$ws.Range("D248").Value = 44558
$ws.Range("D249").Value = 44558
$ws.Range("D250").Value = 44586
$ws.Range("J250").Value = 3100
$ws.Range("D251").Value = 44586
$ws.Range("J251").Value = 1520
$ws.Range("D252").Value = 44294
$ws.Range("J252").Value = 3300
$ws.Range("D253").Value = 44294
$ws.Range("J253").Value = 1540
$ws.Range("D254").Value = 44537
$ws.Range("J254").Value = 3000
$ws.Range("K254").Value = 450
$ws.Range("L254").Value = 500
$ws.Range("M254").Value = 475
$ws.Range("P254").Value = 158
$ws.Range("D255").Value = 44537
$ws.Range("J255").Value = 1500
$ws.Range("K255").Value = 350
$ws.Range("L255").Value = 400
$ws.Range("M255").Value = 375
$ws.Range("P255").Value = 125
$ws.Range("D256").Value = 44824
$ws.Range("J256").Value = 2000
$ws.Range("K256").Value = 550
$ws.Range("L256").Value = 600
$ws.Range("M256").Value = 575
$ws.Range("P256").Value = 192
$ws.Range("D257").Value = 44824
$ws.Range("J257").Value = 1460
$ws.Range("K257").Value = 450
$ws.Range("L257").Value = 500
$ws.Range("M257").Value = 475
$ws.Range("P257").Value = 158
$ws.Range("D258").Value = 44581
$ws.Range("J258").Value = 3080
$ws.Range("K258").Value = 450
$ws.Range("L258").Value = 500
$ws.Range("M258").Value = 475
$ws.Range("P258").Value = 158
$ws.Range("D259").Value = 44581
$ws.Range("J259").Value = 1520
$ws.Range("K259").Value = 350
$ws.Range("L259").Value = 400
$ws.Range("M259").Value = 375
$ws.Range("P259").Value = 125
$ws.Range("D260").Value = 44810
$ws.Range("J260").Value = 2000
$ws.Range("K260").Value = 650
$ws.Range("L260").Value = 700
$ws.Range("M260").Value = 675
$ws.Range("P260").Value = 225
$ws.Range("D261").Value = 44810
$ws.Range("J261").Value = 1480
$ws.Range("K261").Value = 550
$ws.Range("L261").Value = 600
$ws.Range("M261").Value = 575
$ws.Range("P261").Value = 192
$ws.Range("D262").Value = 44203
$ws.Range("K262").Value = 400
$ws.Range("M262").Value = 450
$ws.Range("P262").Value = 150
$ws.Range("D263").Value = 44203
$ws.Range("K263").Value = 300
$ws.Range("L263").Value = 350
$ws.Range("M263").Value = 325
$ws.Range("P263").Value = 108
$ws.Range("D264").Value = 44588
$ws.Range("J264").Value = 2800
$ws.Range("D265").Value = 44588
$ws.Range("J265").Value = 1600
$ws.Range("D266").Value = 44504
$ws.Range("J266").Value = 3000
$ws.Range("D267").Value = 44504
$ws.Range("J267").Value = 1460
$ws.Range("D268").Value = 44670
$ws.Range("J268").Value = 2400
$ws.Range("K268").Value = 450
$ws.Range("M268").Value = 475
$ws.Range("P268").Value = 158
$ws.Range("D269").Value = 44670
$ws.Range("J269").Value = 1500
$ws.Range("K269").Value = 350
$ws.Range("L269").Value = 400
$ws.Range("M269").Value = 375
$ws.Range("P269").Value = 125
$ws.Range("D270").Value = 44187
$ws.Range("J270").Value = 2700
$ws.Range("K270").Value = 400
$ws.Range("L270").Value = 500
$ws.Range("M270").Value = 450
$ws.Range("P270").Value = 150
$ws.Range("D271").Value = 44187
$ws.Range("J271").Value = 1600
$ws.Range("K271").Value = 300
$ws.Range("L271").Value = 350
$ws.Range("M271").Value = 325
$ws.Range("P271").Value = 108
$ws.Range("D272").Value = 44747
$ws.Range("J272").Value = 2440
$ws.Range("K272").Value = 500
$ws.Range("L272").Value = 600
$ws.Range("M272").Value = 550
$ws.Range("P272").Value = 183
$ws.Range("D273").Value = 44747
$ws.Range("J273").Value = 1560
$ws.Range("K273").Value = 400
$ws.Range("L273").Value = 450
$ws.Range("M273").Value = 425
$ws.Range("P273").Value = 142
$ws.Range("D274").Value = 44250
$ws.Range("J274").Value = 3200
$ws.Range("K274").Value = 450
$ws.Range("M274").Value = 475
$ws.Range("P274").Value = 158
$ws.Range("D275").Value = 44250
$ws.Range("J275").Value = 1680
$ws.Range("K275").Value = 350
$ws.Range("L275").Value = 400
$ws.Range("M275").Value = 375
$ws.Range("P275").Value = 125
$ws.Range("D276").Value = 44166
$ws.Range("J276").Value = 2700
$ws.Range("K276").Value = 400
$ws.Range("L276").Value = 500
$ws.Range("M276").Value = 450
$ws.Range("P276").Value = 150
$ws.Range("D277").Value = 44166
$ws.Range("J277").Value = 1600
$ws.Range("K277").Value = 300
$ws.Range("L277").Value = 350
$ws.Range("M277").Value = 325
$ws.Range("P277").Value = 108
$ws.Range("D278").Value = 44635
$ws.Range("J278").Value = 2400
$ws.Range("D279").Value = 44635
$ws.Range("J279").Value = 1300
$ws.Range("D280").Value = 44334
$ws.Range("J280").Value = 3460
$ws.Range("K280").Value = 500
$ws.Range("L280").Value = 600
$ws.Range("M280").Value = 550
$ws.Range("P280").Value = 183
$ws.Range("D281").Value = 44334
$ws.Range("J281").Value = 1680
$ws.Range("K281").Value = 400
$ws.Range("L281").Value = 450
$ws.Range("M281").Value = 425
$ws.Range("P281").Value = 142
$ws.Range("D282").Value = 44427
$ws.Range("J282").Value = 3400
$ws.Range("K282").Value = 450
$ws.Range("L282").Value = 500
$ws.Range("M282").Value = 475
$ws.Range("P282").Value = 158
$ws.Range("D283").Value = 44427
$ws.Range("J283").Value = 1600
$ws.Range("K283").Value = 350
$ws.Range("L283").Value = 400
$ws.Range("M283").Value = 375
$ws.Range("P283").Value = 125
$ws.Range("D284").Value = 44705
$ws.Range("J284").Value = 2400
$ws.Range("K284").Value = 500
$ws.Range("L284").Value = 600
$ws.Range("M284").Value = 550
$ws.Range("P284").Value = 183
$ws.Range("D285").Value = 44705
$ws.Range("J285").Value = 1500
$ws.Range("K285").Value = 400
$ws.Range("L285").Value = 450
$ws.Range("M285").Value = 425
$ws.Range("P285").Value = 142
$ws.Range("D286").Value = 44201
$ws.Range("J286").Value = 2700
$ws.Range("K286").Value = 400
$ws.Range("M286").Value = 450
$ws.Range("P286").Value = 150
$ws.Range("D287").Value = 44201
$ws.Range("J287").Value = 1540
$ws.Range("K287").Value = 300
$ws.Range("L287").Value = 350
$ws.Range("M287").Value = 325
$ws.Range("P287").Value = 108
$ws.Range("D288").Value = 44672
$ws.Range("J288").Value = 2200
$ws.Range("D289").Value = 44672
$ws.Range("J289").Value = 1480
$ws.Range("D290").Value = 44539
$ws.Range("D291").Value = 44539
$ws.Range("J291").Value = 1400
$ws.Range("D292").Value = 44518
$ws.Range("J292").Value = 3000
$ws.Range("D293").Value = 44518
$ws.Range("J293").Value = 1460
$ws.Range("D294").Value = 44488
$ws.Range("J294").Value = 3200
$ws.Range("K294").Value = 450
$ws.Range("L294").Value = 500
$ws.Range("M294").Value = 475
$ws.Range("P294").Value = 158
$ws.Range("D295").Value = 44488
$ws.Range("J295").Value = 1400
$ws.Range("K295").Value = 350
$ws.Range("L295").Value = 400
$ws.Range("M295").Value = 375
$ws.Range("P295").Value = 125
$ws.Range("D296").Value = 44348
$ws.Range("J296").Value = 3400
$ws.Range("D297").Value = 44348
$ws.Range("J297").Value = 1680
$ws.Range("D298").Value = 44614
$ws.Range("J298").Value = 2300
$ws.Range("D299").Value = 44614
$ws.Range("J299").Value = 1460
$ws.Range("D300").Value = 44341
$ws.Range("J300").Value = 3400
$ws.Range("K300").Value = 500
$ws.Range("L300").Value = 600
$ws.Range("M300").Value = 550
$ws.Range("P300").Value = 183
$ws.Range("D301").Value = 44341
$ws.Range("J301").Value = 1680
$ws.Range("K301").Value = 400
$ws.Range("L301").Value = 450
$ws.Range("M301").Value = 425
$ws.Range("P301").Value = 142
$ws.Range("D302").Value = 44278
$ws.Range("J302").Value = 3000
$ws.Range("K302").Value = 450
$ws.Range("M302").Value = 475
$ws.Range("P302").Value = 158
$ws.Range("D303").Value = 44278
$ws.Range("K303").Value = 350
$ws.Range("L303").Value = 400
$ws.Range("M303").Value = 375
$ws.Range("P303").Value = 125
$ws.Range("D304").Value = 44168
$ws.Range("J304").Value = 2800
$ws.Range("K304").Value = 400
$ws.Range("L304").Value = 500
$ws.Range("M304").Value = 450
$ws.Range("P304").Value = 150
$ws.Range("D305").Value = 44168
$ws.Range("J305").Value = 1600
$ws.Range("K305").Value = 300
$ws.Range("L305").Value = 350
$ws.Range("M305").Value = 325
$ws.Range("P305").Value = 108
$ws.Range("D306").Value = 44803
$ws.Range("J306").Value = 2000
$ws.Range("K306").Value = 600
$ws.Range("L306").Value = 700
$ws.Range("M306").Value = 650
$ws.Range("P306").Value = 217
$ws.Range("D307").Value = 44803
$ws.Range("J307").Value = 1400
$ws.Range("K307").Value = 500
$ws.Range("L307").Value = 550
$ws.Range("M307").Value = 525
$ws.Range("P307").Value = 175
$ws.Range("D308").Value = 44222
$ws.Range("J308").Value = 2800
$ws.Range("D309").Value = 44222
$ws.Range("D310").Value = 44224
$ws.Range("D311").Value = 44224
$ws.Range("J311").Value = 1600
$ws.Range("D312").Value = 44462
$ws.Range("J312").Value = 3000
$ws.Range("K312").Value = 450
$ws.Range("L312").Value = 500
$ws.Range("M312").Value = 475
$ws.Range("P312").Value = 158
$ws.Range("D313").Value = 44462
$ws.Range("J313").Value = 1400
$ws.Range("K313").Value = 350
$ws.Range("L313").Value = 400
$ws.Range("M313").Value = 375
$ws.Range("P313").Value = 125
$ws.Range("D314").Value = 44714
$ws.Range("J314").Value = 2480
$ws.Range("K314").Value = 500
$ws.Range("L314").Value = 600
$ws.Range("M314").Value = 550
$ws.Range("P314").Value = 183
$ws.Range("D315").Value = 44714
$ws.Range("J315").Value = 1560
$ws.Range("K315").Value = 400
$ws.Range("L315").Value = 450
$ws.Range("M315").Value = 425
$ws.Range("P315").Value = 142
$ws.Range("D316").Value = 44455
$ws.Range("J316").Value = 3400
$ws.Range("K316").Value = 450
$ws.Range("L316").Value = 500
$ws.Range("M316").Value = 475
$ws.Range("P316").Value = 158
$ws.Range("D317").Value = 44455
$ws.Range("J317").Value = 1600
$ws.Range("K317").Value = 350
$ws.Range("L317").Value = 400
$ws.Range("M317").Value = 375
$ws.Range("P317").Value = 125
$ws.Range("D318").Value = 44357
$ws.Range("J318").Value = 3520
$ws.Range("D319").Value = 44357
$ws.Range("J319").Value = 1680
$ws.Range("D320").Value = 44329
$ws.Range("J320").Value = 3500
$ws.Range("K320").Value = 500
$ws.Range("L320").Value = 600
$ws.Range("M320").Value = 550
$ws.Range("P320").Value = 183
$ws.Range("D321").Value = 44329
$ws.Range("J321").Value = 1660
$ws.Range("K321").Value = 400
$ws.Range("L321").Value = 450
$ws.Range("M321").Value = 425
$ws.Range("P321").Value = 142
$ws.Range("D322").Value = 44495
$ws.Range("J322").Value = 2800
$ws.Range("D323").Value = 44495
$ws.Range("J323").Value = 1400
$ws.Range("D324").Value = 44434
$ws.Range("J324").Value = 3400
$ws.Range("K324").Value = 450
$ws.Range("L324").Value = 500
$ws.Range("M324").Value = 475
$ws.Range("P324").Value = 158
$ws.Range("D325").Value = 44434
$ws.Range("J325").Value = 1600
$ws.Range("K325").Value = 350
$ws.Range("L325").Value = 400
$ws.Range("M325").Value = 375
$ws.Range("P325").Value = 125
$ws.Range("D326").Value = 44628
$ws.Range("J326").Value = 2400
$ws.Range("K326").Value = 500
$ws.Range("L326").Value = 600
$ws.Range("M326").Value = 550
$ws.Range("P326").Value = 183
$ws.Range("D327").Value = 44628
$ws.Range("J327").Value = 1500
$ws.Range("K327").Value = 400
$ws.Range("L327").Value = 450
$ws.Range("M327").Value = 425
$ws.Range("P327").Value = 142
$ws.Range("D328").Value = 44441
$ws.Range("J328").Value = 3300
$ws.Range("D329").Value = 44441
$ws.Range("J329").Value = 1600
$ws.Range("D330").Value = 44273
$ws.Range("J330").Value = 3200
$ws.Range("K330").Value = 450
$ws.Range("L330").Value = 500
$ws.Range("M330").Value = 475
$ws.Range("P330").Value = 158
$ws.Range("D331").Value = 44273
$ws.Range("K331").Value = 350
$ws.Range("L331").Value = 400
$ws.Range("M331").Value = 375
$ws.Range("P331").Value = 125
$ws.Range("D332").Value = 44775
$ws.Range("J332").Value = 2400
$ws.Range("K332").Value = 600
$ws.Range("L332").Value = 700
$ws.Range("M332").Value = 650
$ws.Range("P332").Value = 217
$ws.Range("D333").Value = 44775
$ws.Range("J333").Value = 1500
$ws.Range("K333").Value = 500
$ws.Range("L333").Value = 550
$ws.Range("M333").Value = 525
$ws.Range("P333").Value = 175
$ws.Range("D334").Value = 44362
$ws.Range("J334").Value = 3400
$ws.Range("D335").Value = 44362
$ws.Range("J335").Value = 1700
$ws.Range("D336").Value = 44642
$ws.Range("D337").Value = 44642
$ws.Range("D338").Value = 44707
$ws.Range("J338").Value = 2400
$ws.Range("K338").Value = 500
$ws.Range("L338").Value = 600
$ws.Range("M338").Value = 550
$ws.Range("P338").Value = 183
$ws.Range("D339").Value = 44707
$ws.Range("J339").Value = 1500
$ws.Range("K339").Value = 400
$ws.Range("L339").Value = 450
$ws.Range("M339").Value = 425
$ws.Range("P339").Value = 142
$ws.Range("D340").Value = 44243
$ws.Range("J340").Value = 3000
$ws.Range("D341").Value = 44243
$ws.Range("J341").Value = 1700
$ws.Range("D342").Value = 44579
$ws.Range("J342").Value = 3120
$ws.Range("D343").Value = 44579
$ws.Range("J343").Value = 1560
$ws.Range("D344").Value = 44490
$ws.Range("J344").Value = 3000
$ws.Range("K344").Value = 450
$ws.Range("L344").Value = 500
$ws.Range("M344").Value = 475
$ws.Range("P344").Value = 158
$ws.Range("D345").Value = 44490
$ws.Range("J345").Value = 1400
$ws.Range("K345").Value = 350
$ws.Range("L345").Value = 400
$ws.Range("M345").Value = 375
$ws.Range("P345").Value = 125
$ws.Range("D346").Value = 44845
$ws.Range("K346").Value = 550
$ws.Range("M346").Value = 575
$ws.Range("P346").Value = 192
$ws.Range("D347").Value = 44845
$ws.Range("J347").Value = 1560
$ws.Range("K347").Value = 450
$ws.Range("L347").Value = 500
$ws.Range("M347").Value = 475
$ws.Range("P347").Value = 158
$ws.Range("D348").Value = 44600
$ws.Range("J348").Value = 2400
$ws.Range("K348").Value = 500
$ws.Range("L348").Value = 600
$ws.Range("M348").Value = 550
$ws.Range("P348").Value = 183
$ws.Range("D349").Value = 44600
$ws.Range("J349").Value = 1500
$ws.Range("K349").Value = 400
$ws.Range("L349").Value = 450
$ws.Range("M349").Value = 425
$ws.Range("P349").Value = 142
$ws.Range("D350").Value = 44497
$ws.Range("J350").Value = 2900
$ws.Range("D351").Value = 44497
$ws.Range("J351").Value = 1440
$ws.Range("D352").Value = 44406
$ws.Range("J352").Value = 3500
$ws.Range("D353").Value = 44406
$ws.Range("J353").Value = 1600
$ws.Range("D354").Value = 44280
$ws.Range("J354").Value = 3200
$ws.Range("D355").Value = 44280

# --- Append two new rows (356, 357) with the data pushed out of the rolling window ---
$ws.Range("A356").Value = 8
$ws.Range("B356").Value = "Terminal La Palmera de La Serena"
$ws.Range("C356").Value = "Coquimbo"
$ws.Range("D356").Value = 44572
$ws.Range("D356").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E356").Value = 4
$ws.Range("F356").Value = 100114014
$ws.Range("G356").Value = "Betarraga"
$ws.Range("H356").Value = "Sin especificar"
$ws.Range("I356").Value = "Primera"
$ws.Range("J356").Value = 3160
$ws.Range("K356").Value = 450
$ws.Range("L356").Value = 500
$ws.Range("M356").Value = 475
$ws.Range("N356").Value = "$/paquete 3 unidades"
$ws.Range("O356").Value = "Provincia del Elquí"
$ws.Range("P356").Value = 158
$ws.Range("Q356").Value = 3
$ws.Range("R356").Value = "Hortaliza"
$ws.Range("A357").Value = 8
$ws.Range("B357").Value = "Terminal La Palmera de La Serena"
$ws.Range("C357").Value = "Coquimbo"
$ws.Range("D357").Value = 44572
$ws.Range("D357").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E357").Value = 4
$ws.Range("F357").Value = 100114014
$ws.Range("G357").Value = "Betarraga"
$ws.Range("H357").Value = "Sin especificar"
$ws.Range("I357").Value = "Segunda"
$ws.Range("J357").Value = 1560
$ws.Range("K357").Value = 350
$ws.Range("L357").Value = 400
$ws.Range("M357").Value = 375
$ws.Range("N357").Value = "$/paquete 3 unidades"
$ws.Range("O357").Value = "Provincia del Elquí"
$ws.Range("P357").Value = 125
$ws.Range("Q357").Value = 3
$ws.Range("R357").Value = "Hortaliza"
